# Edit script: add grain quality spec text to column H (rows 41-83)
# and update the sheet view (scroll position / selection) to match the
# author's final state ("OK, ya funciona todo hasta el 6").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('H41').Value = 'Maiz:'
$ws.Range('H42').Value = 'Ph:75kg/hl(mínimo)'
$ws.Range('H43').Value = 'Humedad:14,5%(máximo)'
$ws.Range('H44').Value = 'Granos picados:3,00%(máximo)'
$ws.Range('H46').Value = 'Sorgo granifero:'
$ws.Range('H47').Value = 'Humedad:15%(máximo)'
$ws.Range('H48').Value = 'Granos dañados:2%(máximo)'
$ws.Range('H50').Value = 'Trigo forrajero:'
$ws.Range('H51').Value = 'Ph mínimo: 72,6kg/hl (mínimo)'
$ws.Range('H52').Value = 'Humedad:14%(máximo)'
$ws.Range('H53').Value = 'Insectos:nada,o%'
$ws.Range('H55').Value = 'Cebada cervecera:'
$ws.Range('H56').Value = 'Humedad:12,5%(máximo)'
$ws.Range('H57').Value = 'Proteína mínima:9,5%'
$ws.Range('H58').Value = 'Proteína máxima:13%'
$ws.Range('H60').Value = 'Girasol:'
$ws.Range('H61').Value = 'Humedad:14,0% (máxima)'
$ws.Range('H62').Value = 'Materias grasas:no entiendo la cifra'
$ws.Range('H63').Value = 'Materias extrañas:3,0% (máximo)'
$ws.Range('H65').Value = 'Alpiste:'
$ws.Range('H66').Value = 'Chamico:2 semillas c/100gr'
$ws.Range('H67').Value = 'Humedad:14% (máximo)'
$ws.Range('H68').Value = 'Cornezuelo:0,1%(máximo)'
$ws.Range('H70').Value = 'Avena:'
$ws.Range('H71').Value = 'Ph mínimo:46 /Kg'
$ws.Range('H72').Value = 'Granos manchados: ligeramente'
$ws.Range('H73').Value = 'Humedad: 14%(máximo)'
$ws.Range('H75').Value = 'Mijo:'
$ws.Range('H76').Value = 'Picado:10%(máximo)'
$ws.Range('H77').Value = 'Humedad:15%(máximo)'
$ws.Range('H78').Value = 'Granos dañados:0,50%(máximo)'
$ws.Range('H80').Value = 'Centeno:'
$ws.Range('H81').Value = 'Ph mínimo:70/Kg'
$ws.Range('H82').Value = 'Cornezuelo:0,05%(máximo)'
$ws.Range('H83').Value = 'Tolerancia de picado:1,0%(máximo)'

# Update view state: scroll the window so row 31 / column E is the
# top-left visible cell, and leave the selection on F59 (matches the
# sheetView/selection recorded in the saved workbook).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 5
$ws.Range("F59").Select()

